$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.558.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.30%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.324.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.79%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.11"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.46"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.09%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.56%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.06%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.43%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.894.25"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.70%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.81%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.90"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.504.98"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.15%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.335.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.93%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.85%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.80"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.48"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.06%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.56"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.41%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.34%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.538"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.91%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.78"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.09%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.82"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.25%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0971"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.19%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.32%  "

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.44"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.21%  "

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.97"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.79%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.95"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.89%  "

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.31"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +11.13%  "

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.57"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.29%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.21%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +11.58%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.05"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.21%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +11.98%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.86"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.11%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.842.13"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0734"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.71%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.57%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.33"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.749"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.59%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.59"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.17%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.43%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.05"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +7.86%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.365.87"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.68%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.33%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.97%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.812"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.80%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "282.41"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.23%  "
